# Replace the comma-decimal text values in column K ("SIZE (cm)") with
# real numeric values using a full stop, matching the author's commit
# "Replaced commas with full stops in size columns".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sizeFixes = @{
    3  = 2.5
    4  = 6.5
    10 = 1.5
    12 = 3.5
    15 = 2.5
    31 = 2.5
    35 = 1.5
    55 = 1.5
    56 = 1.5
    57 = 1.5
    58 = 2.5
    59 = 1.5
    62 = 3.5
    63 = 3.5
    70 = 4.5
    71 = 3.5
    72 = 1.5
}

foreach ($row in $sizeFixes.Keys) {
    $ws.Cells.Item($row, 11).Value = $sizeFixes[$row]
}

# Mirror the saved selection state: column K selected (as if the user
# clicked the "SIZE (cm)" column header while editing these values).
$ws.Columns.Item(11).Select()
